$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.932.45'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '1.875.32'
$ws.Range('E3').Value = '  -0.97%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '0.7432'
$ws.Range('E5').Value = '  -3.70%  '
$ws.Range('D6').Value = '242.39'
$ws.Range('E6').Value = '  -0.83%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.3161'
$ws.Range('E8').Value = '  +1.18%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = '0.07188'
$ws.Range('E9').Value = '  -0.64%  '
$ws.Range('B10').Value = 'Solana'
$ws.Range('C10').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').Value = '24.78'
$ws.Range('E10').Value = '  -3.57%  '
$ws.Range('D11').Value = '0.08480'
$ws.Range('E11').Value = '  -3.02%  '
$ws.Range('D12').Value = '0.7546'
$ws.Range('E12').Value = '  -2.17%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.885.22'
$ws.Range('E13').Value = '  -5.88%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '5.398'
$ws.Range('E14').Value = '  -0.52%  '
$ws.Range('D15').Value = '92.78'
$ws.Range('E15').Value = '  -1.61%  '
$ws.Range('D16').Value = '29.935.30'
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('D17').Value = '6.110'
$ws.Range('E17').Value = '  -1.46%  '
$ws.Range('D18').Value = '13.63'
$ws.Range('E18').Value = '  -2.12%  '
$ws.Range('D19').Value = '243.70'
$ws.Range('E19').Value = '  -0.58%  '
$ws.Range('D20').Value = '0.000007832'
$ws.Range('E20').Value = '  -0.58%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').Value = '2.119.26'
$ws.Range('E22').Value = '  -5.19%  '
$ws.Range('D23').Value = '7.990'
$ws.Range('E23').Value = '  -2.22%  '
$ws.Range('D24').Value = '0.9995'
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('D25').Value = '0.1561'
$ws.Range('E25').Value = '  -2.62%  '
$ws.Range('D26').Value = '9.325'
$ws.Range('E26').Value = '  -2.05%  '
$ws.Range('D27').Value = '165.12'
$ws.Range('E27').Value = '  +1.72%  '
$ws.Range('D28').Value = '18.65'
$ws.Range('E28').Value = '  -0.82%  '
$ws.Range('D29').Value = '2.040'
$ws.Range('E29').Value = '  +0.13%  '
$ws.Range('D30').Value = '1.477'
$ws.Range('E30').Value = '  +3.35%  '
$ws.Range('D31').Value = '4.605'
$ws.Range('E31').Value = '  +1.54%  '
$ws.Range('D32').Value = '1.532'
$ws.Range('E32').Value = '  -0.71%  '
$ws.Range('D33').Value = '4.282'
$ws.Range('E33').Value = '  +3.90%  '
$ws.Range('D34').Value = '0.05329'
$ws.Range('E34').Value = '  -2.81%  '
$ws.Range('D35').Value = '1.239'
$ws.Range('E35').Value = '  -0.73%  '
$ws.Range('D36').Value = '0.7569'
$ws.Range('E36').Value = '  +0.50%  '
$ws.Range('D37').Value = '0.9973'
$ws.Range('E37').Value = '  -0.53%  '
$ws.Range('D38').Value = '2.693'
$ws.Range('E38').Value = '  -0.97%  '
$ws.Range('D39').Value = '0.01953'
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('D40').Value = '2.753'
$ws.Range('E40').Value = '  -1.25%  '
$ws.Range('D41').Value = '0.4486'
$ws.Range('E41').Value = '  -0.43%  '
$ws.Range('D42').Value = '1.115.50'
$ws.Range('E42').Value = '  +2.17%  '
$ws.Range('D43').Value = '6.115'
$ws.Range('E43').Value = '  +1.30%  '
$ws.Range('D44').Value = '72.60'
$ws.Range('E44').Value = '  -1.84%  '
$ws.Range('D45').Value = '0.8628'
$ws.Range('E45').Value = '  +0.79%  '
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('D47').Value = '103.09'
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('D48').Value = '7.680'
$ws.Range('E48').Value = '  +0.76%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '1.845'
$ws.Range('E49').Value = '  -2.15%  '
$ws.Range('B50').Value = 'SynthetixNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D50').Value = '3.073'
$ws.Range('E50').Value = '  +4.20%  '
$ws.Range('D51').Value = '2.018.59'
$ws.Range('E51').Value = '  -5.75%  '
